$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H15").Value = 1392.8684
$ws.Range("I15").Value = 1392.8684
$ws.Range("K15").Value = 4178.6052
$ws.Range("M15").Value = -4009.6052

$ws.Range("H38").Value = 800.1539
$ws.Range("I38").Value = 116.833336
$ws.Range("J38").Value = 9000
$ws.Range("K38").Value = 350.500008
$ws.Range("L38").Value = 27000
$ws.Range("M38").Value = 21.49999200000002
$ws.Range("N38").Value = -27744

$ws.Range("H51").Value = 6232.2256
$ws.Range("J51").Value = 5089.478
$ws.Range("L51").Value = 5089.478
$ws.Range("N51").Value = -6057.478

$ws.Range("H61").Value = 357.5
$ws.Range("I61").Value = 357.5
$ws.Range("K61").Value = 1072.5
$ws.Range("M61").Value = -900.5

$ws.Range("H62").Value = 3774.875
$ws.Range("I62").Value = 3069.8333
$ws.Range("J62").Value = 5890
$ws.Range("K62").Value = 3069.8333
$ws.Range("L62").Value = 5890
$ws.Range("M62").Value = -2445.8333
$ws.Range("N62").Value = -7138

$ws.Range("H65").Value = 3774.875
$ws.Range("I65").Value = 3069.8333
$ws.Range("J65").Value = 5890
$ws.Range("K65").Value = 15349.1665
$ws.Range("L65").Value = 29450
$ws.Range("M65").Value = -12229.1665
$ws.Range("N65").Value = -35690

$ws.Range("H109").Value = 96496
$ws.Range("J109").Value = 96496
$ws.Range("L109").Value = 96496
$ws.Range("N109").Value = -99270

$ws.Range("H115").Value = 0
$ws.Range("I115").Value = 0
$ws.Range("K115").Value = 0
$ws.Range("M115").ClearContents()

$ws.Range("H137").Value = 1989.7037
$ws.Range("I137").Value = 1907.826
$ws.Range("K137").Value = 5723.478
$ws.Range("M137").Value = -3173.478

$ws.Range("H138").Value = 5891.528
$ws.Range("J138").Value = 7232.3184
$ws.Range("L138").Value = 21696.9552
$ws.Range("N138").Value = -31976.9552

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 3980.8696
$ws.Range("I2").Value = 4229.3335
$ws.Range("J2").Value = 3515
$ws.Range("K2").Value = 4229.3335
$ws.Range("L2").Value = 3515
$ws.Range("M2").Value = -4116.3335
$ws.Range("N2").Value = -3741

$ws.Range("H32").Value = 1411109.8
$ws.Range("I32").Value = 1517001.4
$ws.Range("K32").Value = 1517001.4
$ws.Range("M32").Value = -1516714.4

$ws.Range("H61").Value = 14377087
$ws.Range("I61").Value = 19092854
$ws.Range("J61").Value = 4002399.2
$ws.Range("K61").Value = 19092854
$ws.Range("L61").Value = 4002399.2
$ws.Range("M61").Value = -19092642
$ws.Range("N61").Value = -4002823.2

$ws.Range("H74").Value = 2100.6365
$ws.Range("I74").Value = 1349
$ws.Range("J74").Value = 5483
$ws.Range("K74").Value = 1349
$ws.Range("L74").Value = 5483
$ws.Range("M74").Value = -475
$ws.Range("N74").Value = -7231

$ws.Range("H77").Value = 2100.6365
$ws.Range("I77").Value = 1349
$ws.Range("J77").Value = 5483
$ws.Range("K77").Value = 6745
$ws.Range("L77").Value = 27415
$ws.Range("M77").Value = -2377
$ws.Range("N77").Value = -36151

$ws.Range("H116").Value = 3980.8696
$ws.Range("I116").Value = 4229.3335
$ws.Range("J116").Value = 3515
$ws.Range("K116").Value = 4229.3335
$ws.Range("L116").Value = 3515
$ws.Range("M116").Value = -1935.3335
$ws.Range("N116").Value = -8103

$ws.Range("H122").Value = 7491.727
$ws.Range("I122").Value = 7756
$ws.Range("K122").Value = 23268
$ws.Range("M122").Value = -20818

$ws.Range("H132").Value = 3131408
$ws.Range("I132").Value = 6271.16
$ws.Range("K132").Value = 18813.48
$ws.Range("M132").Value = -16283.48

$ws.Range("H136").Value = 14377087
$ws.Range("I136").Value = 19092854
$ws.Range("J136").Value = 4002399.2
$ws.Range("K136").Value = 57278562
$ws.Range("L136").Value = 12007197.6
$ws.Range("M136").Value = -57276012
$ws.Range("N136").Value = -12012297.6

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 3980.8696
$ws.Range("I3").Value = 4229.3335
$ws.Range("J3").Value = 3515
$ws.Range("K3").Value = 4229.3335
$ws.Range("L3").Value = 3515
$ws.Range("M3").Value = -4115.3335
$ws.Range("N3").Value = -3743

$ws.Range("H86").Value = 4816.8125
$ws.Range("I86").Value = 1906.4445
$ws.Range("K86").Value = 1906.4445
$ws.Range("M86").Value = -783.4445000000001

$ws.Range("H89").Value = 4816.8125
$ws.Range("I89").Value = 1906.4445
$ws.Range("K89").Value = 9532.2225
$ws.Range("M89").Value = -3916.2225

$ws.Range("H94").Value = 2387.862
$ws.Range("I94").Value = 2401.3572
$ws.Range("K94").Value = 2401.3572
$ws.Range("M94").Value = -1950.3572

$ws.Range("H117").Value = 0
$ws.Range("J117").Value = 0
$ws.Range("L117").ClearContents()
$ws.Range("N117").Value = 0

$ws.Range("H134").Value = 4764604
$ws.Range("I134").Value = 2705.111
$ws.Range("J134").Value = 33335996
$ws.Range("K134").Value = 8115.333
$ws.Range("L134").Value = 100007988
$ws.Range("M134").Value = -5580.333
$ws.Range("N134").Value = -100013058

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H22").Value = 3231.182
$ws.Range("I22").Value = 3236.1428
$ws.Range("J22").Value = 3222.5
$ws.Range("K22").Value = 3236.1428
$ws.Range("L22").Value = 3222.5
$ws.Range("M22").Value = -2886.1428
$ws.Range("N22").Value = -3922.5

$ws.Range("H105").Value = 1994
$ws.Range("I105").Value = 1448.1875
$ws.Range("J105").Value = 3449.5
$ws.Range("K105").Value = 1448.1875
$ws.Range("L105").Value = 3449.5
$ws.Range("M105").Value = 298.8125
$ws.Range("N105").Value = -6943.5

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H4").Value = 2500388.5
$ws.Range("I4").Value = 518
$ws.Range("K4").Value = 1554
$ws.Range("M4").Value = -1442

$ws.Range("H12").Value = 2874
$ws.Range("I12").Value = 201.5
$ws.Range("K12").Value = 604.5
$ws.Range("M12").Value = -431.5

$ws.Range("H17").Value = 12033
$ws.Range("J17").Value = 12033
$ws.Range("L17").Value = 36099
$ws.Range("N17").Value = -36437

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H7").Value = 1000
$ws.Range("J7").Value = 1000
$ws.Range("L7").Value = 1000
$ws.Range("N7").Value = -1224

$ws.Range("H8").Value = 1000
$ws.Range("J8").Value = 1000
$ws.Range("L8").Value = 1000
$ws.Range("N8").Value = -1278

$ws.Range("H29").Value = 10000000
$ws.Range("J29").Value = 0
$ws.Range("L29").Value = 0
$ws.Range("N29").ClearContents()

$ws.Range("H70").Value = 5321.846
$ws.Range("I70").Value = 4825.1816
$ws.Range("K70").Value = 4825.1816
$ws.Range("M70").Value = -4555.1816

$ws.Range("H73").Value = 5321.846
$ws.Range("I73").Value = 4825.1816
$ws.Range("K73").Value = 4825.1816
$ws.Range("M73").Value = -3889.1816

$ws.Range("H113").Value = 2060723.5
$ws.Range("I113").Value = 3499
$ws.Range("J113").Value = 6175172.5
$ws.Range("K113").Value = 3499
$ws.Range("L113").Value = 6175172.5
$ws.Range("M113").Value = -1329
$ws.Range("N113").Value = -6179512.5

$ws.Range("H122").Value = 2299.6924
$ws.Range("I122").Value = 2263.2727
$ws.Range("J122").Value = 2500
$ws.Range("K122").Value = 6789.8181
$ws.Range("L122").Value = 7500
$ws.Range("M122").Value = -4339.8181
$ws.Range("N122").Value = -12400

$ws.Range("H126").Value = 6226.8335
$ws.Range("I126").Value = 6169.5
$ws.Range("K126").Value = 18508.5
$ws.Range("M126").Value = -16038.5

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H42").Value = 74999.5
$ws.Range("J42").Value = 99999
$ws.Range("L42").Value = 99999
$ws.Range("N42").Value = -101125

$ws.Range("H46").Value = 1983.3334
$ws.Range("I46").Value = 1500
$ws.Range("J46").Value = 2950
$ws.Range("K46").Value = 1500
$ws.Range("L46").Value = 2950
$ws.Range("M46").Value = -1312
$ws.Range("N46").Value = -3326

$ws.Range("H48").Value = 700
$ws.Range("I48").Value = 700
$ws.Range("J48").Value = 0
$ws.Range("K48").Value = 700
$ws.Range("L48").ClearContents()
$ws.Range("M48").Value = -39
$ws.Range("N48").Value = 0

$ws.Range("H49").Value = 74999.5
$ws.Range("J49").Value = 99999
$ws.Range("L49").Value = 99999
$ws.Range("N49").Value = -100293

$ws.Range("H93").Value = 3091236.5
$ws.Range("I93").Value = 2256.182
$ws.Range("K93").Value = 2256.182
$ws.Range("M93").Value = -1008.182

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H4").Value = 8499.333000000001
$ws.Range("I4").Value = 11999
$ws.Range("J4").Value = 1500
$ws.Range("K4").Value = 11999
$ws.Range("L4").Value = 1500
$ws.Range("M4").Value = -11886
$ws.Range("N4").Value = -1726

$ws.Range("H54").Value = 33000
$ws.Range("J54").Value = 33000
$ws.Range("L54").Value = 33000
$ws.Range("N54").Value = -34040

$ws.Range("H96").Value = 7134.364
$ws.Range("J96").Value = 5381.3335
$ws.Range("L96").Value = 5381.3335
$ws.Range("N96").Value = -8127.3335

$ws.Range("H122").Value = 2914.0527
$ws.Range("I122").Value = 2804.25
$ws.Range("J122").Value = 3499.6667
$ws.Range("K122").Value = 8412.75
$ws.Range("L122").Value = 10499.0001
$ws.Range("M122").Value = -5962.75
$ws.Range("N122").Value = -15399.0001

$ws.Range("H132").Value = 596756.4399999999
$ws.Range("I132").Value = 8657.267
$ws.Range("J132").Value = 5007500
$ws.Range("K132").Value = 25971.801
$ws.Range("L132").Value = 15022500
$ws.Range("M132").Value = -23441.801
$ws.Range("N132").Value = -15027560

$ws.Range("H136").Value = 1031187.1
$ws.Range("I136").Value = 38553
$ws.Range("J136").Value = 3347333.2
$ws.Range("K136").Value = 115659
$ws.Range("L136").Value = 10041999.6
$ws.Range("M136").Value = -113109
$ws.Range("N136").Value = -10047099.6
